$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells that receive numeric-looking text stay as text (matches original inline-string storage)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '63.724.50'
$ws.Range("E2").Value = '  +1.34%  '

# Row 3
$ws.Range("D3").Value = '3.316.50'
$ws.Range("E3").Value = '  +4.77%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '604.02'
$ws.Range("E5").Value = '  +2.72%  '

# Row 6
$ws.Range("D6").Value = '142.36'
$ws.Range("E6").Value = '  +2.81%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '3.313.24'
$ws.Range("E8").Value = '  +4.84%  '

# Row 9
$ws.Range("D9").Value = '0.519'
$ws.Range("E9").Value = '  +0.62%  '

# Row 10
$ws.Range("E10").Value = '  +2.82%  '

# Row 11
$ws.Range("D11").Value = '5.53'
$ws.Range("E11").Value = '  +4.16%  '

# Row 13
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").Value = '  +1.36%  '

# Row 14
$ws.Range("D14").Value = '34.85'
$ws.Range("E14").Value = '  +2.53%  '

# Row 15
$ws.Range("D15").Value = '3.864.40'
$ws.Range("E15").Value = '  +4.87%  '

# Row 16
$ws.Range("E16").Value = '  -0.01%  '

# Row 17
$ws.Range("D17").Value = '3.315.96'
$ws.Range("E17").Value = '  +4.77%  '

# Row 18
$ws.Range("D18").Value = '63.810.17'
$ws.Range("E18").Value = '  +1.50%  '

# Row 19
$ws.Range("D19").Value = '6.87'
$ws.Range("E19").Value = '  +3.18%  '

# Row 20
$ws.Range("D20").Value = '480.13'
$ws.Range("E20").Value = '  +1.86%  '

# Row 21
$ws.Range("D21").Value = '14.08'
$ws.Range("E21").Value = '  +1.25%  '

# Row 22
$ws.Range("D22").Value = '0.736'
$ws.Range("E22").Value = '  +4.96%  '

# Row 23
$ws.Range("D23").Value = '7.93'
$ws.Range("E23").Value = '  +2.67%  '

# Row 24
$ws.Range("D24").Value = '13.71'
$ws.Range("E24").Value = '  +5.64%  '

# Row 25
$ws.Range("D25").Value = '84.69'
$ws.Range("E25").Value = '  +1.43%  '

# Row 26
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  +2.64%  '

# Row 28
$ws.Range("E28").Value = '  -0.10%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '8.21'
$ws.Range("E29").Value = '  +3.12%  '

# Row 30
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  +2.83%  '

# Row 31
$ws.Range("E31").Value = '  +3.91%  '

# Row 32
$ws.Range("D32").Value = '28.89'
$ws.Range("E32").Value = '  +8.10%  '

# Row 33
$ws.Range("D33").Value = '0.107'
$ws.Range("E33").Value = '  -0.37%  '

# Row 34
$ws.Range("E34").Value = '  +0.45%  '

# Row 35
$ws.Range("E35").Value = '  +3.32%  '

# Row 36
$ws.Range("E36").Value = '  +5.27%  '

# Row 37
$ws.Range("D37").Value = '52.35'
$ws.Range("E37").Value = '  -0.33%  '

# Row 39
$ws.Range("D39").Value = '0.0401'
$ws.Range("E39").Value = '  +3.81%  '

# Row 40
$ws.Range("D40").Value = '434.46'
$ws.Range("E40").Value = '  +4.34%  '

# Row 41
$ws.Range("D41").Value = '3.099.63'
$ws.Range("E41").Value = '  +5.05%  '

# Row 42
$ws.Range("D42").Value = '0.120'
$ws.Range("E42").Value = '  +8.51%  '

# Row 43
$ws.Range("E43").Value = '  +1.35%  '

# Row 44
$ws.Range("E44").Value = '  +0.56%  '

# Row 45
$ws.Range("D45").Value = '0.266'
$ws.Range("E45").Value = '  +1.43%  '

# Row 46
$ws.Range("D46").Value = '2.25'
$ws.Range("E46").Value = '  +6.14%  '

# Row 47
$ws.Range("D47").Value = '37.13'
$ws.Range("E47").Value = '  +16.18%  '

# Row 48
$ws.Range("D48").Value = '26.33'
$ws.Range("E48").Value = '  +3.46%  '

# Row 49
$ws.Range("E49").Value = '  -0.06%  '

# Row 50
$ws.Range("E50").Value = '  +2.70%  '

# Row 51
$ws.Range("E51").Value = '  +0.28%  '
